$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table occupies rows 2..212 (row 1 = header). A new weekly record is
# being inserted at row 103, so every existing record from row 103 downwards
# shifts down by one row, and the record that used to be the very last one
# (row 212) becomes the new last row (213).
#
# Columns A,B,C,E,F,G,H,R are constant for all of these rows, so only the
# "variable" columns D,I,J,K,L,M,N,O,P,Q need to be shifted.
$cols = 4,9,10,11,12,13,14,15,16,17   # D, I, J, K, L, M, N, O, P, Q
$constCols = 1,2,3,5,6,7,8,18         # A, B, C, E, F, G, H, R

# 1) Create the brand new last row (213) as a copy of the current last row
#    (212), before anything else is touched.
foreach ($c in $constCols) {
    $ws.Cells.Item(213, $c).Value = $ws.Cells.Item(212, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(213, $c).Value = $ws.Cells.Item(212, $c).Value2
}
# Keep the date formatting of column D consistent with the rest of the table.
$ws.Cells.Item(213, 4).NumberFormat = $ws.Cells.Item(212, 4).NumberFormat

# 2) Shift rows 212 down to 104 so that each row takes on the values that
#    belonged to the row above it. Walking from the bottom upwards guarantees
#    that the "source" row for each step still holds its original value when
#    it is read.
for ($r = 212; $r -ge 104; $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($r - 1, $c).Value2
    }
}

# 3) Finally, write the new incoming record into row 103.
$ws.Cells.Item(103, 4).Value = 44494
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 80
$ws.Cells.Item(103, 11).Value = 7000
$ws.Cells.Item(103, 12).Value = 7000
$ws.Cells.Item(103, 13).Value = 7000
$ws.Cells.Item(103, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(103, 15).Value = "Región del Maule"
$ws.Cells.Item(103, 16).Value = 350
$ws.Cells.Item(103, 17).Value = 20
